# Auto-generated edit script: updates Leve market-price / profit columns (H:N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets to reflect a refreshed
# market-data snapshot ("chore: update Sheets via scheduled runner").
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 2310
$ws.Range("I58").Value = 1500
$ws.Range("J58").Value = 2425.7144
$ws.Range("K58").Value = 4500
$ws.Range("L58").Value = 7277.1432
$ws.Range("M58").Value = -4350
$ws.Range("N58").Value = -7577.1432
$ws.Range("H113").Value = 3362.7
$ws.Range("I113").Value = 3361.0715
$ws.Range("K113").Value = 3361.0715
$ws.Range("M113").Value = -107.0715
$ws.Range("H126").Value = 49980
$ws.Range("J126").Value = 49980
$ws.Range("L126").Value = 49980
$ws.Range("N126").Value = -59860
$ws.Range("H132").Value = 2864.3333
$ws.Range("I132").Value = 2985.7646
$ws.Range("J132").Value = 800
$ws.Range("K132").Value = 8957.293799999999
$ws.Range("L132").Value = 2400
$ws.Range("M132").Value = -6427.293799999999
$ws.Range("N132").Value = -7460
$ws.Range("H137").Value = 959.53845
$ws.Range("I137").Value = 897.63635
$ws.Range("J137").Value = 1300
$ws.Range("K137").Value = 2692.90905
$ws.Range("L137").Value = 3900
$ws.Range("M137").Value = -142.9090500000002
$ws.Range("N137").Value = -9000

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3026.7896
$ws.Range("I2").Value = 2649.8333
$ws.Range("J2").Value = 3673
$ws.Range("K2").Value = 2649.8333
$ws.Range("L2").Value = 3673
$ws.Range("M2").Value = -2536.8333
$ws.Range("N2").Value = -3899
$ws.Range("H116").Value = 3026.7896
$ws.Range("I116").Value = 2649.8333
$ws.Range("J116").Value = 3673
$ws.Range("K116").Value = 2649.8333
$ws.Range("L116").Value = 3673
$ws.Range("M116").Value = -355.8332999999998
$ws.Range("N116").Value = -8261
$ws.Range("H132").Value = 1363.6571
$ws.Range("I132").Value = 883.1667
$ws.Range("K132").Value = 2649.5001
$ws.Range("M132").Value = -119.5001000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3026.7896
$ws.Range("I3").Value = 2649.8333
$ws.Range("J3").Value = 3673
$ws.Range("K3").Value = 2649.8333
$ws.Range("L3").Value = 3673
$ws.Range("M3").Value = -2535.8333
$ws.Range("N3").Value = -3901
$ws.Range("H36").Value = 6559.8
$ws.Range("I36").Value = 2150
$ws.Range("J36").Value = 7049.778
$ws.Range("K36").Value = 2150
$ws.Range("L36").Value = 7049.778
$ws.Range("M36").Value = -1616
$ws.Range("N36").Value = -8117.778
$ws.Range("H37").Value = 8412
$ws.Range("J37").Value = 8412
$ws.Range("L37").Value = 8412
$ws.Range("N37").Value = -8686
$ws.Range("H39").Value = 6052.75
$ws.Range("J39").Value = 6052.75
$ws.Range("L39").Value = 6052.75
$ws.Range("N39").Value = -6830.75
$ws.Range("H105").Value = 3234.9092
$ws.Range("I105").Value = 2490.3215
$ws.Range("J105").Value = 7404.6
$ws.Range("K105").Value = 2490.3215
$ws.Range("L105").Value = 7404.6
$ws.Range("M105").Value = -743.3215
$ws.Range("N105").Value = -10898.6
$ws.Range("H134").Value = 39647.965
$ws.Range("I134").Value = 2076
$ws.Range("J134").Value = 171149.83
$ws.Range("K134").Value = 6228
$ws.Range("L134").Value = 513449.49
$ws.Range("M134").Value = -3693
$ws.Range("N134").Value = -518519.49

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 883.1667
$ws.Range("J15").Value = 1266.3334
$ws.Range("L15").Value = 1266.3334
$ws.Range("N15").Value = -1606.3334
$ws.Range("H118").Value = 30000
$ws.Range("J118").Value = 30000
$ws.Range("L118").Value = 30000
$ws.Range("N118").Value = -33314
$ws.Range("H134").Value = 38462860
$ws.Range("I134").Value = 1443.5454
$ws.Range("J134").Value = 250000660
$ws.Range("K134").Value = 4330.6362
$ws.Range("L134").Value = 750001980
$ws.Range("M134").Value = -1795.6362
$ws.Range("N134").Value = -750007050

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 782.5
$ws.Range("I33").Value = 575
$ws.Range("J33").Value = 990
$ws.Range("K33").Value = 3450
$ws.Range("L33").Value = 5940
$ws.Range("M33").Value = -3167
$ws.Range("N33").Value = -6506
$ws.Range("H64").Value = 1335.5
$ws.Range("J64").Value = 1390.5714
$ws.Range("L64").Value = 4171.7142
$ws.Range("N64").Value = -4711.7142
$ws.Range("H67").Value = 1335.5
$ws.Range("J67").Value = 1390.5714
$ws.Range("L67").Value = 4171.7142
$ws.Range("N67").Value = -6043.7142
$ws.Range("H96").Value = 3532.611
$ws.Range("J96").Value = 3532.611
$ws.Range("L96").Value = 10597.833
$ws.Range("N96").Value = -14715.833
$ws.Range("H113").Value = 29140.656
$ws.Range("I113").Value = 1250
$ws.Range("J113").Value = 30831
$ws.Range("K113").Value = 3750
$ws.Range("L113").Value = 92493
$ws.Range("M113").Value = -1580
$ws.Range("N113").Value = -96833

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 5000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -5540
$ws.Range("H73").Value = 5000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 5000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -6872
$ws.Range("H107").Value = 440.33334
$ws.Range("I107").Value = 289.44446
$ws.Range("J107").Value = 666.6667
$ws.Range("K107").Value = 289.44446
$ws.Range("L107").Value = 666.6667
$ws.Range("M107").Value = 1630.55554
$ws.Range("N107").Value = -4506.6667
$ws.Range("H122").Value = 941519.4399999999
$ws.Range("I122").Value = 1317397.2
$ws.Range("J122").Value = 1825
$ws.Range("K122").Value = 3952191.6
$ws.Range("L122").Value = 1825
$ws.Range("M122").Value = -3949741.6
$ws.Range("N122").Value = -10375
$ws.Range("H132").Value = 2838.4243
$ws.Range("I132").Value = 2688.238
$ws.Range("J132").Value = 3101.25
$ws.Range("K132").Value = 8064.714
$ws.Range("L132").Value = 9303.75
$ws.Range("M132").Value = -5534.714
$ws.Range("N132").Value = -14363.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 433.33334
$ws.Range("H27").Value = 433.33334
$ws.Range("H50").Value = 6166.6665
$ws.Range("J50").Value = 6166.6665
$ws.Range("L50").Value = 6166.6665
$ws.Range("N50").Value = -7440.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 17520
$ws.Range("I32").Value = 10000
$ws.Range("J32").Value = 19400
$ws.Range("K32").Value = 10000
$ws.Range("L32").Value = 19400
$ws.Range("M32").Value = -9683
$ws.Range("N32").Value = -20034
$ws.Range("H39").Value = 8000
$ws.Range("I39").Value = 8000
$ws.Range("K39").Value = 8000
$ws.Range("M39").Value = -7587
